$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number format on cells whose new values would otherwise be
# auto-converted to numeric values by Excel (single decimal point numbers),
# so they remain plain text matching the original inline-string cell type.
$textCells = @("D5", "D6", "D8", "D9", "D10", "D11", "D12", "D14", "D19", "D20", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D31", "D32", "D33", "D34", "D36", "D38", "D39", "D41", "D42", "D43", "D46", "D47", "D50", "D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '48.199.02'
$ws.Range("E2").Value = '  +2.13%  '
$ws.Range("D3").Value = '2.526.68'
$ws.Range("E3").Value = '  +1.27%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '324.08'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("D6").Value = '109.13'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  +0.88%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").Value = '0.560'
$ws.Range("E9").Value = '  +4.77%  '
$ws.Range("D10").Value = '40.77'
$ws.Range("E10").Value = '  +4.39%  '
$ws.Range("D11").Value = '20.35'
$ws.Range("E11").Value = '  +10.94%  '
$ws.Range("D12").Value = '0.0829'
$ws.Range("E12").Value = '  +1.97%  '
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("D14").Value = '7.29'
$ws.Range("E14").Value = '  +1.45%  '
$ws.Range("D15").Value = '2.919.69'
$ws.Range("E15").Value = '  +0.98%  '
$ws.Range("D16").Value = '2.531.47'
$ws.Range("E16").Value = '  +1.18%  '
$ws.Range("E17").Value = '  +1.13%  '
$ws.Range("D18").Value = '48.031.20'
$ws.Range("E18").Value = '  +1.85%  '
$ws.Range("D19").Value = '13.27'
$ws.Range("E19").Value = '  +3.54%  '
$ws.Range("D20").Value = '6.65'
$ws.Range("E20").Value = '  +0.36%  '
$ws.Range("D21").Value = '0.0₃0956'
$ws.Range("E21").Value = '  +1.63%  '
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("D23").Value = '72.37'
$ws.Range("E23").Value = '  +2.34%  '
$ws.Range("D24").Value = '270.25'
$ws.Range("E24").Value = '  +9.48%  '
$ws.Range("D25").Value = '2.57'
$ws.Range("E25").Value = '  -0.72%  '
$ws.Range("D26").Value = '26.27'
$ws.Range("E26").Value = '  +1.03%  '
$ws.Range("E27").Value = '  -0.27%  '
$ws.Range("D28").Value = '10.17'
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("D29").Value = '0.146'
$ws.Range("E29").Value = '  +4.82%  '
$ws.Range("D30").Value = '35.80'
$ws.Range("E30").Value = '  +1.73%  '
$ws.Range("D31").Value = '2.11'
$ws.Range("E31").Value = '  -8.23%  '
$ws.Range("D32").Value = '49.75'
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("D33").Value = '19.99'
$ws.Range("E33").Value = '  +0.05%  '
$ws.Range("D34").Value = '5.42'
$ws.Range("E34").Value = '  -0.02%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").Value = '0.0796'
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("E37").Value = '  +1.39%  '
$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  +1.36%  '
$ws.Range("D39").Value = '3.00'
$ws.Range("E39").Value = '  +0.79%  '
$ws.Range("E40").Value = '  +0.13%  '
$ws.Range("D41").Value = '22.28'
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").Value = '119.30'
$ws.Range("E42").Value = '  -2.62%  '
$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").Value = '2.19'
$ws.Range("E43").Value = '  -2.08%  '
$ws.Range("E44").Value = '  +0.59%  '
$ws.Range("D45").Value = '2.011.71'
$ws.Range("E45").Value = '  +1.09%  '
$ws.Range("D46").Value = '3.14'
$ws.Range("E46").Value = '  +2.78%  '
$ws.Range("D47").Value = '2.05'
$ws.Range("E47").Value = '  -0.04%  '
$ws.Range("E48").Value = '  +5.47%  '
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").Value = '5.26'
$ws.Range("E50").Value = '  +1.42%  '
$ws.Range("D51").Value = '79.76'
$ws.Range("E51").Value = '  +2.40%  '
